# Adds extra columns (WIN, TOP4, TOP5, RELEGATION) ahead of ExpPoints,
# which now moves from column C to column G. ExpPoints values are
# refreshed, and a few teams swap positions in the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# The original header style (bold font, thin border, centered) lived
# only on A1:C1. Copy it across the newly introduced D1:G1 header
# cells too, then fill in the new header labels.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("G1").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# --- New team order + refreshed ExpPoints (now column G) --------------
# rowIndex -> (Team, ExpPoints)
$teams = @{
    2  = @("Arsenal", 79.40248878607713)
    3  = @("Liverpool", 71.26751588073536)
    4  = @("Manchester City", 70.88678636532417)
    5  = @("Chelsea", 62.19724749174479)
    6  = @("Aston Villa", 57.08971354027942)
    7  = @("Crystal Palace", 56.86135978515252)
    8  = @("Newcastle United", 56.04147495373084)
    9  = @("Tottenham Hotspur", 55.38666095710316)
    10 = @("Brighton & Hove Albion", 52.54613251582245)
    11 = @("AFC Bournemouth", 52.17430150585)
    12 = @("Brentford", 51.45918798040534)
    13 = @("Manchester United", 50.10506306229097)
    14 = @("Everton", 46.04335435571845)
    15 = @("Fulham", 43.32597134712751)
    16 = @("Nottingham Forest", 42.75216109960289)
    17 = @("Sunderland", 37.97983027614948)
    18 = @("Leeds United", 36.26630914223142)
    19 = @("West Ham United", 36.26170765736258)
    20 = @("Burnley", 35.95772203182079)
    21 = @("Wolverhampton Wanderers", 31.70763586628734)
}

foreach ($row in $teams.Keys) {
    $team = $teams[$row][0]
    $expPoints = $teams[$row][1]

    $ws.Cells.Item($row, 2).Value = $team

    # Clear out the newly inserted WIN/TOP4/TOP5/RELEGATION cells (C-F)
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""

    # ExpPoints now lives in column G
    $ws.Cells.Item($row, 7).Value = $expPoints
}
